# Add a simple engine with GameObject, Components and Prefabs support:
# three new sheets -> DayNightCycle, Animals, PrefabsView

$wb = $excel.ActiveWorkbook
$tileTypes = $wb.Worksheets.Item("Tile Types")

# xlPasteFormats
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Sheet 2: DayNightCycle
# ---------------------------------------------------------------------
$dayNight = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $tileTypes)
$dayNight.Name = "DayNightCycle"

# header row
$dayNight.Range("A1").Value = "IDS"
$dayNight.Range("A1").Font.ThemeColor = 1
$dayNight.Range("B1").Value = "Name"
$dayNight.Range("C1").Value = "Duration"

# reuse the newly-created style on the rest of the sheet via copy/paste
$dayNight.Range("A1").Copy()
$dayNight.Range("B1:C1").PasteSpecial($xlPasteFormats)
$dayNight.Range("A2:C7").PasteSpecial($xlPasteFormats)

$dayParts = @("Morning", "MidDay", "Afternoon", "Evening", "Sunrise", "Dusk")
for ($i = 0; $i -lt $dayParts.Length; $i++) {
    $r = $i + 2
    $dayNight.Cells.Item($r, 1).Value = $dayParts[$i]
    $dayNight.Cells.Item($r, 2).Value = $dayParts[$i]
    $dayNight.Cells.Item($r, 3).Value = 60
}

# ---------------------------------------------------------------------
# Sheet 3: Animals
# ---------------------------------------------------------------------
$animals = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dayNight)
$animals.Name = "Animals"

$animals.Range("A1").Value = "IDS"
$animals.Range("B1").Value = "Name"
$animals.Range("A2").Value = "Monkey"
$animals.Range("B2").Value = "Monkey"

# style 1 (same as Tile Types header/data cells) for A1:B2
$tileTypes.Range("A2").Copy()
$animals.Range("A1:B2").PasteSpecial($xlPasteFormats)

# style 3 (same as Tile Types blank padding rows) for the rest of the grid
$tileTypes.Range("A9").Copy()
$animals.Range("C1:K2").PasteSpecial($xlPasteFormats)
$animals.Range("A3:K24").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Sheet 4: PrefabsView
# ---------------------------------------------------------------------
$prefabsView = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $animals)
$prefabsView.Name = "PrefabsView"

$prefabsView.Range("A1").Value = "IDS"
$prefabsView.Range("A1").Font.ThemeColor = 1
$prefabsView.Range("B1").Value = "Architecture ID"
$prefabsView.Range("C1").Value = "Prefab resource path"
$prefabsView.Range("A2").Value = "Monkey view"
$prefabsView.Range("B2").Value = "Monkey"
$prefabsView.Range("C2").Value = "../Prefabs/Monkey.xml"

$prefabsView.Range("A1").Copy()
$prefabsView.Range("B1:C1").PasteSpecial($xlPasteFormats)
$prefabsView.Range("A2:C2").PasteSpecial($xlPasteFormats)

$prefabsView.Columns.Item(1).ColumnWidth = 10.833333333333334
$prefabsView.Columns.Item(2).ColumnWidth = 13.166666666666666
$prefabsView.Columns.Item(3).ColumnWidth = 37.416666666666664

$tileTypes.Activate()
